# "rename TTML pre-processor to TTML Segmenter"
#
# The shape named "Rectangle 24" (inside group "Groupe 13") holds a
# two-line label "TTML " / "Post-Processor" built from several runs
# separated by a soft line break (<a:br/>). We only want to retarget
# the word that changed ("Post-Processor" -> "Segmenter") while leaving
# the paragraph/run/line-break structure alone, so we locate the
# sub-range of characters that spells "Post-Processor" and overwrite
# just that text - this keeps the existing <a:br/> intact instead of
# collapsing the paragraph into a single run or splitting it in two.

$oldWord = "Post-Processor"
$newWord = "Segmenter"

function Update-ShapeText($shape) {
    if ($shape.Type -eq 6) {
        # msoGroup: recurse into the group's members.
        $items = $shape.GroupItems
        for ($i = 1; $i -le $items.Count; $i++) {
            Update-ShapeText $items.Item($i)
        }
        return
    }

    if (-not $shape.HasTextFrame) {
        return
    }
    if (-not $shape.TextFrame.HasText) {
        return
    }

    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($oldWord)
    while ($idx -ge 0) {
        $sub = $tr.Characters($idx + 1, $oldWord.Length)
        $sub.Text = $newWord
        $full = $tr.Text
        $idx = $full.IndexOf($oldWord)
    }
}

$p = $ppt.ActivePresentation
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        Update-ShapeText $slide.Shapes.Item($shi)
    }
}
